# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# This script operates on the "ODI Batting Extra" worksheet (sheet index 3)
# of the workbook. 15 new MATCH_CODE records were scraped that chronologically
# precede the previously-known earliest record, so they are inserted as new
# rows at the top of the data table (rows 2-16), pushing the existing 20 data
# rows down by 15 (to rows 17-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Insert 15 blank rows above the current row 2, shifting all existing data
# (previously rows 2-21) down to rows 17-36.
$ws.Range("A2:A16").EntireRow.Insert()

# New rows of data scraped for the extra batting fields, keyed by row number.
# Each entry gives: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL
# Rows with no batting data still get a MAN_OF_MATCH value of "NO".
$newData = @{
    2  = @{ A = "3006"; B = 7;    C = "5"; D = "0"; E = "8.41%" }
    3  = @{ A = "3009"; B = $null; C = $null; D = $null; E = $null }
    4  = @{ A = "3010"; B = $null; C = $null; D = $null; E = $null }
    5  = @{ A = "3012"; B = $null; C = $null; D = $null; E = $null }
    6  = @{ A = "3015"; B = 2;    C = "5"; D = "0"; E = "23.08%" }
    7  = @{ A = "3017"; B = 2;    C = "3"; D = "0"; E = "5.30%" }
    8  = @{ A = "3018"; B = $null; C = $null; D = $null; E = $null }
    9  = @{ A = "3019"; B = 2;    C = "1"; D = "0"; E = "2.27%" }
    10 = @{ A = "3024"; B = $null; C = $null; D = $null; E = $null }
    11 = @{ A = "3028"; B = $null; C = $null; D = $null; E = $null }
    12 = @{ A = "3030"; B = $null; C = $null; D = $null; E = $null }
    13 = @{ A = "3032"; B = 2;    C = "1"; D = "0"; E = "1.55%" }
    14 = @{ A = "3034"; B = $null; C = $null; D = $null; E = $null }
    15 = @{ A = "3040"; B = $null; C = $null; D = $null; E = $null }
    16 = @{ A = "3042"; B = 2;    C = "1"; D = "0"; E = "3.14%" }
}

foreach ($r in 2..16) {
    $row = $newData[$r]

    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row.A

    if ($null -eq $row.B) {
        $ws.Range("B$r").Text = ""
    } else {
        $ws.Range("B$r").Value = $row.B
    }

    if ($null -eq $row.C) {
        $ws.Range("C$r").Text = ""
    } else {
        $ws.Range("C$r").NumberFormat = "@"
        $ws.Range("C$r").Value = $row.C
    }

    if ($null -eq $row.D) {
        $ws.Range("D$r").Text = ""
    } else {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $row.D
    }

    if ($null -eq $row.E) {
        $ws.Range("E$r").Text = ""
    } else {
        $ws.Range("E$r").NumberFormat = "@"
        $ws.Range("E$r").Value = $row.E
    }

    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = "NO"
}
